# Update "want-to-go" counts (column F) for a few rows on both the
# "展览" sheet and the "全部类型" sheet, reflecting refreshed scrape data.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 657
$wsExpo.Range("F4").Value = 1514
$wsExpo.Range("F5").Value = 703
$wsExpo.Range("F6").Value = 10

# 全部类型 (All types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 657
$wsAll.Range("F4").Value = 1514
$wsAll.Range("F6").Value = 703
$wsAll.Range("F7").Value = 10
